$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card14")

$rng = $ws.Range("A3:A12")
$rng.NumberFormat = "@"
for ($r = 3; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "14"
}
